$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.646.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '''3.444.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''580.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = '''148.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("E9").Value = '  +4.32%  '
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("E11").Value = '  +2.20%  '
$ws.Range("D12").Value = '''4.036.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("E13").Value = '  +1.98%  '
$ws.Range("E14").Value = '  -5.63%  '
$ws.Range("D15").Value = '''3.451.67'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '''62.711.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '''14.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("E20").Value = '  -3.08%  '
$ws.Range("D21").Value = '''386.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").Value = '''0.563'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '''3.581.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.21%  '
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("E30").Value = '  -2.94%  '
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -4.60%  '
$ws.Range("D34").Value = '''23.24'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("D35").Value = '''1.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.50%  '
$ws.Range("E36").Value = '  +1.63%  '
$ws.Range("D37").Value = '''31.89'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("D39").Value = '''169.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("E41").Value = '  +0.62%  '
$ws.Range("E42").Value = '  -2.51%  '
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("D47").Value = '''2.570.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("D50").Value = '''22.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.00%  '
$ws.Range("E51").Value = '  -0.06%  '
